# Generate Report for Handback
# Update the timestamp cells that record when the handback report / xliff
# generation occurred, reflecting a regenerated report a bit later than before.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 15:25:58"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 15:25:53"
$wsZhCn.Range("K2").Value = "2016-08-30 15:26:25"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 15:25:58"
$wsDeDe.Range("K2").Value = "2016-08-30 15:26:32"
